# Updates the cryptos price-tracker sheet (Sheet1) with refreshed
# Price (D) / Volume(1h) (E) figures, matching a scheduled GitHub
# Actions data refresh. A new coin ("ApeXProtocol") was inserted at
# row 42, pushing the former rows 42-50 down to 43-51 and dropping the
# previous last row ("Stacks").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '42.932.18'
$ws.Range('E2').Value = '  -1.31%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.342.80'
$ws.Range('E3').Value = '  +1.35%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.01%  '

# Row 5: BNB
$ws.Range('E5').Value = '  -1.36%  '

# Row 6: Solana
$ws.Range('D6').Value = '''101.21'
$ws.Range('E6').Value = '  -0.91%  '

# Row 7: XRP
$ws.Range('D7').Value = '''0.510'
$ws.Range('E7').Value = '  -4.96%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  -3.42%  '

# Row 10: Avalanche
$ws.Range('D10').Value = '''34.96'
$ws.Range('E10').Value = '  -2.33%  '

# Row 11: OKB
$ws.Range('E11').Value = '  +0.08%  '

# Row 12: Dogecoin
$ws.Range('D12').Value = '''0.0802'
$ws.Range('E12').Value = '  -1.87%  '

# Row 13: TRON
$ws.Range('E13').Value = '  -0.29%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''6.80'
$ws.Range('E14').Value = '  -3.07%  '

# Row 15: Chainlink
$ws.Range('D15').Value = '''15.88'
$ws.Range('E15').Value = '  +6.12%  '

# Row 16: WrappedEther
$ws.Range('D16').Value = '2.279.44'
$ws.Range('E16').Value = '  -1.41%  '

# Row 17: Polygon
$ws.Range('E17').Value = '  +0.38%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '42.860.08'
$ws.Range('E18').Value = '  -1.26%  '

# Row 19: Uniswap
$ws.Range('D19').Value = '''6.26'
$ws.Range('E19').Value = '  +1.38%  '

# Row 20: ShibaInu
$ws.Range('E20').Value = '  -1.69%  '

# Row 21: InternetComputer(DFINITY)
$ws.Range('D21').Value = '''11.69'
$ws.Range('E21').Value = '  -5.99%  '

# Row 22: Litecoin
$ws.Range('D22').Value = '''67.89'
$ws.Range('E22').Value = '  -0.45%  '

# Row 23: BitcoinCash
$ws.Range('D23').Value = '''236.85'
$ws.Range('E23').Value = '  -2.14%  '

# Row 24: ImmutableX
$ws.Range('D24').Value = '''2.04'
$ws.Range('E24').Value = '  +0.18%  '

# Row 25: PancakeSwap
$ws.Range('D25').Value = '''2.57'
$ws.Range('E25').Value = '  -2.06%  '

# Row 26: Dai
$ws.Range('E26').Value = '  -0.04%  '

# Row 27: EthereumClassic
$ws.Range('D27').Value = '''25.50'
$ws.Range('E27').Value = '  +2.85%  '

# Row 28: Toncoin
$ws.Range('E28').Value = '  +1.07%  '

# Row 29: InjectiveProtocol
$ws.Range('D29').Value = '''35.11'
$ws.Range('E29').Value = '  -4.24%  '

# Row 30: Cosmos
$ws.Range('D30').Value = '''9.36'
$ws.Range('E30').Value = '  -2.89%  '

# Row 31: Monero
$ws.Range('D31').Value = '''160.53'
$ws.Range('E31').Value = '  -4.28%  '

# Row 32: FirstDigitalUSD
$ws.Range('D32').Value = '''1.00'

# Row 33: Filecoin
$ws.Range('E33').Value = '  -2.88%  '

# Row 34: RenderToken
$ws.Range('D34').Value = '''4.67'
$ws.Range('E34').Value = '  +7.82%  '

# Row 35: WEMIXToken
$ws.Range('E35').Value = '  -0.29%  '

# Row 36: Celestia
$ws.Range('D36').Value = '''17.44'
$ws.Range('E36').Value = '  -0.75%  '

# Row 37: Hedera
$ws.Range('E37').Value = '  -1.96%  '

# Row 38: LidoDAOToken
$ws.Range('D38').Value = '''2.97'
$ws.Range('E38').Value = '  -4.01%  '

# Row 39: ARBITRUM
$ws.Range('E39').Value = '  -0.14%  '

# Row 40: Kaspa
$ws.Range('E40').Value = '  -3.17%  '

# Row 41: Stellar
$ws.Range('E41').Value = '  -2.41%  '

# Row 42: ApeXProtocol
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').Value = '''2.47'
$ws.Range('E42').Value = '  +7.00%  '

# Row 43: Maker
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.022.25'
$ws.Range('E43').Value = '  +2.61%  '

# Row 44: VeChain
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0286'
$ws.Range('E44').Value = '  -1.26%  '

# Row 45: EnergySwap
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''18.76'
$ws.Range('E45').Value = '  -3.84%  '

# Row 46: FraxShare
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '''10.31'
$ws.Range('E46').Value = '  +3.36%  '

# Row 47: NEARProtocol
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '''2.94'
$ws.Range('E47').Value = '  -1.35%  '

# Row 48: MultiversX
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').Value = '''56.39'
$ws.Range('E48').Value = '  +1.42%  '

# Row 49: HuobiToken
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').Value = '''2.90'
$ws.Range('E49').Value = '  -1.04%  '

# Row 50: RocketPoolETH
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.567.61'
$ws.Range('E50').Value = '  +1.14%  '

# Row 51: THORChain
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '''4.66'
$ws.Range('E51').Value = '  +1.68%  '
